# Update scripts with new TPM values (Fgf5-Fgfr3 LR-pair sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Delete the trailing rows (old rows 5, 6, 7) ---------------------------
# These correspond to target clusters MuSCs/Neutrophils/Resolving-Mac that no
# longer appear in the refreshed TPM output.
$ws.Rows("5:7").Delete() | Out-Null

# --- Row 4: target cluster changes from Inflammatory-Mac to MuSCs ----------
$ws.Range("D4").Value = "MuSCs"

# --- Row 2 numeric updates ---------------------------------------------------
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.902492
$ws.Range("H2").Value = 3.804984
$ws.Range("M2").Value = 7.6704545
$ws.Range("N2").Value = 15.340909
$ws.Range("O2").Value = 0.8278663930876066
$ws.Range("P2").Value = 0.7913005936208135
$ws.Range("Q2").Value = 14.592978322614
$ws.Range("R2").Value = 58.371913290456
$ws.Range("S2").Value = 0.8278663930876066
$ws.Range("T2").Value = 0.7913005936208135

# --- Row 3 numeric updates ---------------------------------------------------
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.902492
$ws.Range("H3").Value = 3.804984
$ws.Range("O3").Value = 0.0924194920655273
$ws.Range("P3").Value = 0.1325061620042962
$ws.Range("Q3").Value = 1.629098192124
$ws.Range("R3").Value = 9.774589152743999
$ws.Range("S3").Value = 0.0924194920655273
$ws.Range("T3").Value = 0.1325061620042962

# --- Row 4 numeric updates ---------------------------------------------------
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.902492
$ws.Range("H4").Value = 3.804984
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.7385775
$ws.Range("N4").Value = 1.477155
$ws.Range("O4").Value = 0.0797141148468662
$ws.Range("P4").Value = 0.07619324437489022
$ws.Range("Q4").Value = 1.40513778513
$ws.Range("R4").Value = 5.62055114052
$ws.Range("S4").Value = 0.0797141148468662
$ws.Range("T4").Value = 0.07619324437489022
